$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2189.7273
$ws.Range("I28").Value = 239.85715
$ws.Range("J28").Value = 5602
$ws.Range("K28").Value = 239.85715
$ws.Range("L28").Value = 5602
$ws.Range("M28").Value = 245.14285
$ws.Range("N28").Value = -6572
$ws.Range("H62").Value = 13462
$ws.Range("I62").Value = 13462
$ws.Range("K62").Value = 13462
$ws.Range("M62").Value = -12838
$ws.Range("H65").Value = 13462
$ws.Range("I65").Value = 13462
$ws.Range("K65").Value = 67310
$ws.Range("M65").Value = -64190
$ws.Range("H76").Value = 3124.5
$ws.Range("I76").Value = 3170.8572
$ws.Range("K76").Value = 3170.8572
$ws.Range("M76").Value = -2855.8572
$ws.Range("H79").Value = 3124.5
$ws.Range("I79").Value = 3170.8572
$ws.Range("K79").Value = 3170.8572
$ws.Range("M79").Value = -2078.8572
$ws.Range("H112").Value = 1526.9849
$ws.Range("J112").Value = 1526.9849
$ws.Range("L112").Value = 4580.9547
$ws.Range("N112").Value = -6796.9547
$ws.Range("H137").Value = 48861.855
$ws.Range("J137").Value = 92408.45
$ws.Range("L137").Value = 277225.35
$ws.Range("N137").Value = -282325.35
$ws.Range("H138").Value = 1531.8
$ws.Range("J138").Value = 1775.8813
$ws.Range("L138").Value = 5327.6439
$ws.Range("N138").Value = -15607.6439

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 4000080
$ws.Range("I8").Value = 4000080
$ws.Range("K8").Value = 4000080
$ws.Range("M8").Value = -3999936
$ws.Range("H32").Value = 5498.5054
$ws.Range("I32").Value = 3726.3171
$ws.Range("J32").Value = 16676.924
$ws.Range("K32").Value = 3726.3171
$ws.Range("L32").Value = 16676.924
$ws.Range("M32").Value = -3439.3171
$ws.Range("N32").Value = -17250.924
$ws.Range("H74").Value = 607.625
$ws.Range("I74").Value = 607.625
$ws.Range("K74").Value = 607.625
$ws.Range("M74").Value = 266.375
$ws.Range("H77").Value = 607.625
$ws.Range("I77").Value = 607.625
$ws.Range("K77").Value = 3038.125
$ws.Range("M77").Value = 1329.875
$ws.Range("H97").Value = 1138.8636
$ws.Range("I97").Value = 1103.3158
$ws.Range("J97").Value = 1364
$ws.Range("K97").Value = 1103.3158
$ws.Range("L97").Value = 1364
$ws.Range("M97").Value = -607.3158000000001
$ws.Range("N97").Value = -2356
$ws.Range("H102").Value = 2775.3333
$ws.Range("I102").Value = 2760.5
$ws.Range("K102").Value = 2760.5
$ws.Range("M102").Value = -1138.5
$ws.Range("H122").Value = 3850.7144
$ws.Range("I122").Value = 1489.75
$ws.Range("K122").Value = 4469.25
$ws.Range("M122").Value = -2019.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4548481
$ws.Range("I20").Value = 6669185.5
$ws.Range("K20").Value = 6669185.5
$ws.Range("M20").Value = -6668938.5
$ws.Range("H134").Value = 3776.9678
$ws.Range("I134").Value = 3615.6667
$ws.Range("K134").Value = 10847.0001
$ws.Range("M134").Value = -8312.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 211
$ws.Range("I7").Value = 243.4
$ws.Range("K7").Value = 243.4
$ws.Range("M7").Value = -130.4
$ws.Range("H31").Value = 3249.4
$ws.Range("I31").Value = 1000
$ws.Range("K31").Value = 1000
$ws.Range("M31").Value = -705
$ws.Range("H34").Value = 3249.4
$ws.Range("I34").Value = 1000
$ws.Range("K34").Value = 1000
$ws.Range("M34").Value = -798
$ws.Range("H58").Value = 1812764.6
$ws.Range("I58").Value = 2558817.5
$ws.Range("J58").Value = 921.8570999999999
$ws.Range("K58").Value = 2558817.5
$ws.Range("L58").Value = 921.8570999999999
$ws.Range("M58").Value = -2558614.5
$ws.Range("N58").Value = -1327.8571
$ws.Range("H122").Value = 3156.5
$ws.Range("I122").Value = 1987.8
$ws.Range("K122").Value = 5963.4
$ws.Range("M122").Value = -3513.4
$ws.Range("H136").Value = 1812764.6
$ws.Range("I136").Value = 2558817.5
$ws.Range("J136").Value = 921.8570999999999
$ws.Range("K136").Value = 7676452.5
$ws.Range("L136").Value = 2765.5713
$ws.Range("M136").Value = -7673902.5
$ws.Range("N136").Value = -7865.5713

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1103.2273
$ws.Range("I122").Value = 598.5
$ws.Range("J122").Value = 1153.7
$ws.Range("K122").Value = 5386.5
$ws.Range("L122").Value = 10383.3
$ws.Range("M122").Value = -2936.5
$ws.Range("N122").Value = -15283.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5945314
$ws.Range("I11").Value = 7699970
$ws.Range("J11").Value = 2436001.5
$ws.Range("K11").Value = 7699970
$ws.Range("L11").Value = 2436001.5
$ws.Range("M11").Value = -7699831
$ws.Range("N11").Value = -2436279.5
$ws.Range("H20").Value = 1403778.9
$ws.Range("I20").Value = 2100001.8
$ws.Range("J20").Value = 11333.333
$ws.Range("K20").Value = 2100001.8
$ws.Range("L20").Value = 11333.333
$ws.Range("M20").Value = -2099756.8
$ws.Range("N20").Value = -11823.333
$ws.Range("H97").Value = 1534.8667
$ws.Range("J97").Value = 1382.6666
$ws.Range("L97").Value = 1382.6666
$ws.Range("N97").Value = -2374.6666
$ws.Range("H102").Value = 2101.111
$ws.Range("I102").Value = 2238.75
$ws.Range("K102").Value = 2238.75
$ws.Range("M102").Value = -616.75
$ws.Range("H136").Value = 17857.5
$ws.Range("J136").Value = 17857.5
$ws.Range("L136").Value = 53572.5
$ws.Range("N136").Value = -58672.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3030.3
$ws.Range("I7").Value = 2800.5
$ws.Range("K7").Value = 2800.5
$ws.Range("M7").Value = -2688.5
$ws.Range("H126").Value = 3030.3
$ws.Range("I126").Value = 2800.5
$ws.Range("K126").Value = 8401.5
$ws.Range("M126").Value = -5931.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 80005
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H21").Value = 14833.333
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 14833.333
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 14833.333
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -15303.333
$ws.Range("H31").Value = 15180
$ws.Range("J31").Value = 15180
$ws.Range("L31").Value = 15180
$ws.Range("N31").Value = -15876
$ws.Range("H35").Value = 14833.333
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 14833.333
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 14833.333
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -15413.333
$ws.Range("H126").Value = 1425.8276
$ws.Range("I126").Value = 1224.4
$ws.Range("K126").Value = 3673.2
$ws.Range("M126").Value = -1203.2
